# 自动更新Excel文件 - daily "remaining days" rollover
#
# For every data row (row 2..last used row) on the active sheet:
#   - Column E ("剩余" / remaining days) is decremented by 1 for the new day.
#   - When a row's remaining-day counter has run out (E == 1), the cycle
#     restarts: E is reset to 10 and the start date in column F ("开始时间",
#     stored as a literal yyyymmdd integer) is pushed forward by 10 days.
#   - Rows whose F value is not a well-formed 8-digit yyyymmdd date (e.g. a
#     corrupted/placeholder value) are left untouched, since the date math
#     can't be performed on them.
#
# Columns: A=行号 B=店铺名称 C=地址 D=总天 E=剩余 F=开始时间 G=备注1 H=备注2 I=备注3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {

    $eCell = $ws.Cells.Item($row, 5)
    $fCell = $ws.Cells.Item($row, 6)

    $eVal = $eCell.Value()
    $fVal = $fCell.Value()

    if ($eVal -eq $null -or $fVal -eq $null) {
        continue
    }

    $fStr = [string]([int64]$fVal)

    # Only touch rows with a well-formed 8-digit yyyymmdd start date.
    if ($fStr.Length -ne 8) {
        continue
    }

    $year = [int]$fStr.Substring(0, 4)
    $month = [int]$fStr.Substring(4, 2)
    $day = [int]$fStr.Substring(6, 2)

    $validDate = $true
    try {
        $startDate = Get-Date -Year $year -Month $month -Day $day
    } catch {
        $validDate = $false
    }

    if (-not $validDate) {
        continue
    }

    $eNum = [int]$eVal

    if ($eNum -eq 1) {
        # Cycle restarts: remaining days reset to 10, start date rolls
        # forward by 10 days.
        $newE = 10
        $newDate = $startDate.AddDays(10)
        $newF = [int]($newDate.ToString("yyyyMMdd"))
    } else {
        # One more day has elapsed; remaining days decreases by 1.
        $newE = $eNum - 1
        $newF = [int64]$fVal
    }

    $eCell.Value = $newE
    $fCell.Value = $newF
}
